# refactor DialogLoader data structure
# Adds a "FinalResults" column (L) and a pair of new jump-target values
# ("Jump1"/"Jump2") on row 11 (columns D/E) of the DialogLoader table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the additional "FinalResults" column.
$ws.Range("L1").Value = "FinalResults"

# New data cells on row 11 for the @Jump entry.
$ws.Range("D11").Value = "Jump1"
$ws.Range("E11").Value = "Jump2"
$ws.Range("L11").Value = "H2"

# Match the author's final cursor position/viewport from the commit.
$ws.Activate()
$ws.Range("N6").Select()
